# A new daily price record (Mango, Vega Central Mapocho de Santiago) was
# inserted ahead of the existing row 634, pushing the former rows 634-720
# down to 635-721 (dimension grows from A1:T720 to A1:T721).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 634; everything below (634-720) shifts down
# one row to (635-721), carrying its data/formatting with it.
$ws.Rows("634:634").Insert()

# Populate the newly inserted row 634 with the new record's data.
$ws.Range("A634").Value = 9
$ws.Range("B634").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C634").Value = "Metropolitana"
$ws.Range("D634").Value = 45142
$ws.Range("E634").Value = 13
$ws.Range("F634").Value = "Fruta"
$ws.Range("G634").Value = 100108
$ws.Range("H634").Value = "Tropicales y subtropicales"
$ws.Range("I634").Value = 100108002
$ws.Range("J634").Value = "Mango"
$ws.Range("K634").Value = "Sin especificar"
$ws.Range("L634").Value = "Primera"
$ws.Range("M634").Value = 590
$ws.Range("N634").Value = 7500
$ws.Range("O634").Value = 8000
$ws.Range("P634").Value = 7746
$ws.Range("Q634").Value = "$/bandeja 4 kilos"
$ws.Range("R634").Value = "Brasil"
$ws.Range("S634").Value = 1936
$ws.Range("T634").Value = 4
